$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3099390012751145
$ws.Range("J2").Value = 0.3099390012751145
$ws.Range("M2").Value = 2.231113333333334
$ws.Range("N2").Value = 6.69334
$ws.Range("O2").Value = 0.01598125358798882
$ws.Range("P2").Value = 0.01598125358798882
$ws.Range("Q2").Value = 0.2675402368444444
$ws.Range("R2").Value = 2.4078621316
$ws.Range("S2").Value = 0.004953213776185595
$ws.Range("T2").Value = 0.004953213776185595
$ws.Range("I3").Value = 0.3099390012751145
$ws.Range("J3").Value = 0.3099390012751145
$ws.Range("O3").Value = 0.1634493267640196
$ws.Range("P3").Value = 0.1634493267640195
$ws.Range("S3").Value = 0.05065932109633008
$ws.Range("T3").Value = 0.05065932109633007
$ws.Range("I4").Value = 0.3099390012751145
$ws.Range("J4").Value = 0.3099390012751145
$ws.Range("M4").Value = 58.02175166666666
$ws.Range("N4").Value = 174.065255
$ws.Range("O4").Value = 0.4156043142904646
$ws.Range("P4").Value = 0.4156043142904646
$ws.Range("Q4").Value = 6.957581648188888
$ws.Range("R4").Value = 62.61823483369999
$ws.Range("S4").Value = 0.1288119860968154
$ws.Range("T4").Value = 0.1288119860968154
$ws.Range("I5").Value = 0.3099390012751145
$ws.Range("J5").Value = 0.3099390012751145
$ws.Range("M5").Value = 15.16934033333333
$ws.Range("N5").Value = 45.508021
$ws.Range("O5").Value = 0.1086565487318021
$ws.Range("P5").Value = 0.1086565487318021
$ws.Range("Q5").Value = 1.819006163837778
$ws.Range("R5").Value = 16.37105547454
$ws.Range("S5").Value = 0.03367690219593556
$ws.Range("T5").Value = 0.03367690219593556
$ws.Range("I6").Value = 0.3099390012751145
$ws.Range("J6").Value = 0.3099390012751145
$ws.Range("M6").Value = 41.36709099999999
$ws.Range("N6").Value = 124.101273
$ws.Range("O6").Value = 0.2963085566257249
$ws.Range("P6").Value = 0.2963085566257249
$ws.Range("Q6").Value = 4.960465772113333
$ws.Range("R6").Value = 44.64419194902
$ws.Range("S6").Value = 0.0918375781098479
$ws.Range("T6").Value = 0.0918375781098479
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2669800000000001
$ws.Range("H7").Value = 0.8009400000000001
$ws.Range("I7").Value = 0.6900609987248855
$ws.Range("J7").Value = 0.6900609987248854
$ws.Range("M7").Value = 2.231113333333334
$ws.Range("N7").Value = 6.69334
$ws.Range("O7").Value = 0.01598125358798882
$ws.Range("P7").Value = 0.01598125358798882
$ws.Range("Q7").Value = 0.5956626377333335
$ws.Range("R7").Value = 5.360963739600001
$ws.Range("S7").Value = 0.01102803981180322
$ws.Range("T7").Value = 0.01102803981180322
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2669800000000001
$ws.Range("H8").Value = 0.8009400000000001
$ws.Range("I8").Value = 0.6900609987248855
$ws.Range("J8").Value = 0.6900609987248854
$ws.Range("O8").Value = 0.1634493267640196
$ws.Range("P8").Value = 0.1634493267640195
$ws.Range("Q8").Value = 6.092178975820001
$ws.Range("R8").Value = 54.82961078238
$ws.Range("S8").Value = 0.1127900056676895
$ws.Range("T8").Value = 0.1127900056676895
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2669800000000001
$ws.Range("H9").Value = 0.8009400000000001
$ws.Range("I9").Value = 0.6900609987248855
$ws.Range("J9").Value = 0.6900609987248854
$ws.Range("M9").Value = 58.02175166666666
$ws.Range("N9").Value = 174.065255
$ws.Range("O9").Value = 0.4156043142904646
$ws.Range("P9").Value = 0.4156043142904646
$ws.Range("Q9").Value = 15.49064725996667
$ws.Range("R9").Value = 139.4158253397
$ws.Range("S9").Value = 0.2867923281936492
$ws.Range("T9").Value = 0.2867923281936491
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2669800000000001
$ws.Range("H10").Value = 0.8009400000000001
$ws.Range("I10").Value = 0.6900609987248855
$ws.Range("J10").Value = 0.6900609987248854
$ws.Range("M10").Value = 15.16934033333333
$ws.Range("N10").Value = 45.508021
$ws.Range("O10").Value = 0.1086565487318021
$ws.Range("P10").Value = 0.1086565487318021
$ws.Range("Q10").Value = 4.049910482193334
$ws.Range("R10").Value = 36.44919433974
$ws.Range("S10").Value = 0.07497964653586656
$ws.Range("T10").Value = 0.07497964653586654
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2669800000000001
$ws.Range("H11").Value = 0.8009400000000001
$ws.Range("I11").Value = 0.6900609987248855
$ws.Range("J11").Value = 0.6900609987248854
$ws.Range("M11").Value = 41.36709099999999
$ws.Range("N11").Value = 124.101273
$ws.Range("O11").Value = 0.2963085566257249
$ws.Range("P11").Value = 0.2963085566257249
$ws.Range("Q11").Value = 11.04418595518
$ws.Range("R11").Value = 99.39767359662001
$ws.Range("S11").Value = 0.204470978515877
$ws.Range("T11").Value = 0.204470978515877

Write-Output "applied $(($ws.Range("A1").Worksheet.Name)) updates"
